# Reposition "Bottom Stringer1" (row 11) to its new y'/z' coordinates.
# D11 and E11 were hard-coded/simple-formula values; the new iteration
# drives them from the stringer-placement constants below (entered as
# formulas, matching how the author recorded them in Excel).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D11").Formula = "=-1.6046754518"
$ws.Range("E11").Formula = "=-0.092597515"
